$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.944.39'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '2.235.38'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  +0.15%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '313.21'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -2.08%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '98.69'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -5.43%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.569'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -3.74%  '
$ws.Range('E8').Value = '  +0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.533'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -7.62%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.19'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -6.65%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0820'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -2.89%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '7.37'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -7.14%  '
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D14').Value = '2.575.86'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '2.237.72'
$ws.Range('E15').Value = '  -2.35%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.837'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -5.15%  '
$ws.Range('E17').Value = '  -3.95%  '
$ws.Range('D18').Value = '43.811.06'
$ws.Range('E18').Value = '  -0.83%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.96'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -12.65%  '
$ws.Range('D20').Value = '0.0₃0963'
$ws.Range('E20').Value = '  -3.84%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.33'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -5.29%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '65.01'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('E23').Value = '  -7.56%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '233.35'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('E25').Value = '  -8.95%  '
$ws.Range('E26').Value = '  +0.39%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.16'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.17%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.77%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '36.68'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -8.36%  '
$ws.Range('E30').Value = '  -9.20%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '157.91'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -2.21%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '19.91'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -3.37%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0830'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -7.07%  '
$ws.Range('E34').Value = '  -1.45%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.20'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -6.64%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.110'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +1.49%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.90'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -6.78%  '
$ws.Range('E38').Value = '  -3.85%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '15.70'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('E40').Value = '  -9.67%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '4.04'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -11.52%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0307'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -6.78%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').Value = '1.709.89'
$ws.Range('E44').Value = '  -5.58%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.194'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -7.74%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '5.14'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -5.66%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '73.14'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -5.09%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '80.07'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -8.03%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.66'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.06%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '101.70'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -2.93%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '56.39'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -6.22%  '
